$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "harvester" column (B) for all data rows (2-24): it incorrectly
# held the same text as rnaSampleNumber ("Retrofitted_0779"); it should
# instead hold the harvester initials "S.GISH".
$ws.Range("B2:B24").Value = "S.GISH"

# Leave the selection on the harvester column, matching where the edit
# was made.
$ws.Range("B:B").Select()
